# --- Dodanie podziału treningu na części (Trening column) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Trening" header, styled the same as the other headers (bold + border).
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

$ws.Cells.Item(2, 1).Value = 45685.65130914352
$ws.Cells.Item(2, 2).Value = 1388.1
$ws.Cells.Item(2, 3).Value = 13.98
$ws.Cells.Item(2, 4).Value = 3.541718551090787
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45685.65796076389
$ws.Cells.Item(3, 2).Value = 1962.8
$ws.Cells.Item(3, 3).Value = 13.63
$ws.Cells.Item(3, 4).Value = 3.666645833424162
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45685.65866215277
$ws.Cells.Item(4, 2).Value = 2023.4
$ws.Cells.Item(4, 3).Value = 14.09
$ws.Cells.Item(4, 4).Value = 3.484031813485283
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45685.64713321759
$ws.Cells.Item(5, 2).Value = 1027.3
$ws.Cells.Item(5, 3).Value = 8.89
$ws.Cells.Item(5, 4).Value = 3.236595051629203
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45685.65130567129
$ws.Cells.Item(6, 2).Value = 1387.8
$ws.Cells.Item(6, 3).Value = 9.75
$ws.Cells.Item(6, 4).Value = 3.267182792936052
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45685.65602673611
$ws.Cells.Item(7, 2).Value = 1795.7
$ws.Cells.Item(7, 3).Value = 9.26
$ws.Cells.Item(7, 4).Value = 3.101976701191493
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45685.66909502315
$ws.Cells.Item(8, 2).Value = 2924.8
$ws.Cells.Item(8, 3).Value = 14.68
$ws.Cells.Item(8, 4).Value = 4.258690834045412
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45685.6777744213
$ws.Cells.Item(9, 2).Value = 3674.7
$ws.Cells.Item(9, 3).Value = 13.22
$ws.Cells.Item(9, 4).Value = 4.141723905290876
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45685.6837443287
$ws.Cells.Item(10, 2).Value = 4190.5
$ws.Cells.Item(10, 3).Value = 14.17
$ws.Cells.Item(10, 4).Value = 3.861694676535471
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45685.66909155092
$ws.Cells.Item(11, 2).Value = 2924.5
$ws.Cells.Item(11, 3).Value = 9.31
$ws.Cells.Item(11, 4).Value = 3.125275343656541
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45685.66951168981
$ws.Cells.Item(12, 2).Value = 2960.8
$ws.Cells.Item(12, 3).Value = 9.01
$ws.Cells.Item(12, 4).Value = 2.941868884222849
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45685.67777210648
$ws.Cells.Item(13, 2).Value = 3674.5
$ws.Cells.Item(13, 3).Value = 8.65
$ws.Cells.Item(13, 4).Value = 3.41021989924567
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# Apply the datetime number format to column A (rows 2-13).
# The first cell goes through both format revisions (164 then 165) like the
# original authoring session; the rest just take the final format directly so
# the style table ends up with exactly one extra cellXf (the one actually used).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
